$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 data: Subject MVI010R141, Ear R, Date 8/25/2021 (serial 44433)
$ws.Range("A11").Value = "MVI010R141"
$ws.Range("B11").Value = "R"

# Copy the date cell format from C10 (numFmtId 14 date style) onto C11, then set value
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C11").Value2 = 44433

# Update selection to reflect new active cell D11
$ws.Range("D11").Select()
